# The deck's single slide master (and therefore every slide) points at
# ppt/theme/theme2.xml ("Simple Light"), while the unused notes-master theme
# part (ppt/theme/theme1.xml) holds the "Default" color scheme. The upstream
# edit swaps the two themes' contents: the deck's active theme becomes the
# "Default" palette (and the idle notes-master theme becomes "Simple Light").
#
# Re-point every slide's live theme color scheme (backed by the shared
# ppt/theme/theme2.xml part) at the "Default" palette's twelve colors, in
# DrawingML clrScheme order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.

$p = $ppt.ActivePresentation

$defaultTheme = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    5800213,    # dk2      158158
    15987699,   # lt2      F3F3F3
    13077765,   # accent1  058DC7
    3322960,    # accent2  50B432
    1791725,    # accent3  ED561B
    61421,      # accent4  EDEF00
    15059748,   # accent5  24CBE5
    7529828,    # accent6  64E572
    13369378,   # hlink    2200CC
    9116245     # folHlink 551A8B
)

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $tcs = $p.Slides.Item($si).ThemeColorScheme
    for ($i = 1; $i -le $tcs.Count; $i++) {
        $tcs.Item($i).RGB = $defaultTheme[$i - 1]
    }
}
